# Applies the "Add files via upload" revision to 01项目计划表.xlsx:
#  1. Marks the six "2018.10.25 第八周周四" tasks (rows 93-98) as 已完成 (column C).
#  2. Appends a brand-new weekly block ("日期：2018.10.29 第九周周一", rows 101-110)
#     built the same way the sheet's earlier weekly blocks are: a merged title
#     row, the 组员/计划内容/完成情况/备注 header row, six member rows, and a
#     merged two-row 总结： summary placeholder.
#  3. The long-form summary paragraph that used to live in the merged
#     A99:D100 cell moves down to become the new block's summary (A109:D110);
#     the old A99:D100 cell is repurposed to hold this revision's own summary
#     paragraph about the API work that was completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in "已完成" for the already-finished tasks of the prior week ---
$ws.Range("C93").Value = "已完成"
$ws.Range("C94").Value = "已完成"
$ws.Range("C95").Value = "已完成"
$ws.Range("C96").Value = "已完成"
$ws.Range("C97").Value = "已完成"
$ws.Range("C98").Value = "已完成"

# --- 2. Move the old summary paragraph text down, put the new one in its place ---
$ws.Range("A99").Value = "总结：经过几天阅读demo和查询基本上已经可以实现api接口，统一返回json数据。由于先前没有撰写api接口文档，所以下周将进行api接口实现及文档撰写，为了与前端能够尽快实现交互，具体为每写一个api接口文档即实现该api，两者同步进行。此外，下周开始，app端进行逻辑实现。"

# --- 3. Build rows 101-110 by cloning the formatting of the previous block (rows 91-100) ---
$ws.Range("A91:D100").Copy()
$ws.Range("A101").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Title row for the new week
$ws.Range("A101").Value = "日期：2018.10.29 第九周周一"

# Header row
$ws.Range("A102").Value = "组员"
$ws.Range("B102").Value = "计划内容"
$ws.Range("C102").Value = "完成情况"
$ws.Range("D102").Value = "备注"

# Member rows
$ws.Range("A103").Value = "李光洪"
$ws.Range("B103").Value = "1、登录api；2、注册api；3、查询用户api"
$ws.Range("D103").Value = "api实现及api接口文档撰写"

$ws.Range("A104").Value = "吴彤林"
$ws.Range("B104").Value = "app我的界面逻辑编码"

$ws.Range("A105").Value = "劳汉文"
$ws.Range("B105").Value = "app通讯录界面逻辑编码"

$ws.Range("A106").Value = "方嘉耀"
$ws.Range("B106").Value = "app地图界面逻辑编码"

$ws.Range("A107").Value = "成世靖"
$ws.Range("B107").Value = "app消息界面逻辑编码"

$ws.Range("A108").Value = "丰浩"
$ws.Range("B108").Value = "协作app通讯录界面逻辑编码"

# Summary placeholder (the text that used to occupy A99:D100)
$ws.Range("A109").Value = "总结："

# --- 4. Re-create the merges for the new block ---
$ws.Range("A101:D101").Merge()
$ws.Range("A109:D110").Merge()

# --- 5. Restore the view state Excel saved the file with ---
$ws.Range("C103").Select()
$excel.ActiveWindow.ScrollRow = 82
